# Temporary Fix for RestDayEquality Penalty Value:
# Separate exhibition games from the list of total games and only use
# non-exhibition games when computing/summing the penalty stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: "Quality" -> "Penalty" for the first (Entire games) table ---
$ws.Range("B1").Value2 = "Penalty"

# --- Division U7 rows (2-4), now excluding exhibition games ---
$ws.Range("B2").Value2 = 509.45
$ws.Range("C2").Value2 = 276.0
$ws.Range("D2").Value2 = 10.0
$ws.Range("E2").Value2 = 266.0
$ws.Range("F2").Value2 = "1 sec"

$ws.Range("B3").Value2 = 709.27
$ws.Range("C3").Value2 = 173.0
$ws.Range("D3").Value2 = 2.0
$ws.Range("E3").Value2 = 171.0
$ws.Range("F3").Value2 = "1 sec"

$ws.Range("B4").Value2 = 851.87
$ws.Range("C4").Value2 = 301.0
$ws.Range("D4").Value2 = 18.0
$ws.Range("E4").Value2 = 283.0
$ws.Range("F4").Value2 = "2 sec"

# Sums row only over the 3 real rows now (2:4) instead of (2:12)
$ws.Range("B13").Formula = "=SUM(B2:B4)"
$ws.Range("C13").Formula = "=SUM(C2:C4)"
$ws.Range("D13").Formula = "=SUM(D2:D4)"
$ws.Range("E13").Formula = "=SUM(E2:E4)"

# --- Division tier-3 rows (16-18) ---
$ws.Range("B16").Value2 = 307.45
$ws.Range("C16").Value2 = 1692.0
$ws.Range("D16").Value2 = 44.0
$ws.Range("E16").Value2 = 1648.0

$ws.Range("B17").Value2 = 454.57
$ws.Range("C17").Value2 = 1035.0
$ws.Range("D17").Value2 = 40.0
$ws.Range("E17").Value2 = 995.0

$ws.Range("B18").Value2 = 741.87
$ws.Range("C18").Value2 = 1812.0
$ws.Range("D18").Value2 = 58.0
$ws.Range("E18").Value2 = 1754.0

# Sums row only over the 3 real rows now (16:18) instead of (16:26)
$ws.Range("B27").Formula = "=SUM(B16:B18)"
$ws.Range("C27").Formula = "=SUM(C16:C18)"
$ws.Range("D27").Formula = "=SUM(D16:D18)"
$ws.Range("E27").Formula = "=SUM(E16:E18)"

# --- Entire League row (30) ---
$ws.Range("B30").Value2 = 1283.02
$ws.Range("C30").Value2 = 9168.0
$ws.Range("D30").Value2 = 223.0
$ws.Range("E30").Value2 = 8945.0
$ws.Range("F30").Value2 = "0 min, 11 sec"
